$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("B2").Value = 14.91337859113037
$ws.Range("C2").Value = 7.212868576558786
$ws.Range("E2").Value = 16.54446560179825
$ws.Range("F2").Value = 47.91076956138986
$ws.Range("G2").Value = 3.709998007382997
$ws.Range("J2").Value = 9.650641237403047
$ws.Range("K2").Value = 14.38396544212531
$ws.Range("N2").Value = 21.89333937772982

# Row 3
$ws.Range("B3").Value = 14.69750152644504
$ws.Range("C3").Value = 7.048716882043077
$ws.Range("E3").Value = 16.34417675229887
$ws.Range("F3").Value = 47.6468337714286
$ws.Range("G3").Value = 3.713076479700827
$ws.Range("J3").Value = 9.672004212554022
$ws.Range("K3").Value = 14.24343735969443
$ws.Range("N3").Value = 21.94678794939214

# Row 4
$ws.Range("B4").Value = 14.56777913382746
$ws.Range("C4").Value = 6.948366528903175
$ws.Range("E4").Value = 16.22489733094899
$ws.Range("F4").Value = 47.49680410918433
$ws.Range("G4").Value = 3.715063444799427
$ws.Range("J4").Value = 9.686331621631279
$ws.Range("K4").Value = 14.16043410670485
$ws.Range("N4").Value = 21.98155345644411

# Row 5
$ws.Range("B5").Value = 14.51569588720611
$ws.Range("C5").Value = 6.907649200289464
$ws.Range("E5").Value = 16.177277846496
$ws.Range("F5").Value = 47.4387331393232
$ws.Range("G5").Value = 3.71589757288675
$ws.Range("J5").Value = 9.692474436548054
$ws.Range("K5").Value = 14.12747026955175
$ws.Range("N5").Value = 21.9962105622945

# Row 6
$ws.Range("B6").Value = 14.50709657600154
$ws.Range("C6").Value = 6.900900709998243
$ws.Range("E6").Value = 16.16943197616614
$ws.Range("F6").Value = 47.42927691011432
$ws.Range("G6").Value = 3.716037556989884
$ws.Range("J6").Value = 9.693512821083555
$ws.Range("K6").Value = 14.12204963241091
$ws.Range("N6").Value = 21.99867395357952

# Row 7
$ws.Range("B7").Value = 14.56707347428178
$ws.Range("C7").Value = 6.947816599325293
$ws.Range("E7").Value = 16.22425104346007
$ws.Range("F7").Value = 47.49600847166251
$ws.Range("G7").Value = 3.71507459513352
$ws.Range("J7").Value = 9.686413233902343
$ws.Range("K7").Value = 14.15998601459908
$ws.Range("N7").Value = 21.98174914375264

# Row 8
$ws.Range("B8").Value = 14.8383991307027
$ws.Range("C8").Value = 7.156216872843594
$ws.Range("E8").Value = 16.47467362709805
$ws.Range("F8").Value = 47.81729528112966
$ws.Range("G8").Value = 3.711039434135619
$ws.Range("J8").Value = 9.657755901799456
$ws.Range("K8").Value = 14.33485210530936
$ws.Range("N8").Value = 21.91136420701605

# Row 9
$ws.Range("B9").Value = 15.38982032614113
$ws.Range("C9").Value = 7.565570908348608
$ws.Range("E9").Value = 16.99246779069043
$ws.Range("F9").Value = 48.54066435362424
$ws.Range("G9").Value = 3.703890219438893
$ws.Range("J9").Value = 9.6111680590927
$ws.Range("K9").Value = 14.70216059448507
$ws.Range("N9").Value = 21.78879128321948

# Row 10
$ws.Range("B10").Value = 15.80259549736125
$ws.Range("C10").Value = 7.863144877769644
$ws.Range("E10").Value = 17.38556141738001
$ws.Range("F10").Value = 49.12588197348525
$ws.Range("G10").Value = 3.699097519282732
$ws.Range("J10").Value = 9.582804652259194
$ws.Range("K10").Value = 14.98460217898417
$ws.Range("N10").Value = 21.70815096868315

# Row 11
$ws.Range("B11").Value = 15.99116502496876
$ws.Range("C11").Value = 7.99712596753743
$ws.Range("E11").Value = 17.56636048707995
$ws.Range("F11").Value = 49.40303063376574
$ws.Range("G11").Value = 3.69701581304721
$ws.Range("J11").Value = 9.571176332673783
$ws.Range("K11").Value = 15.11531670214195
$ws.Range("N11").Value = 21.67350903194389

# Row 12
$ws.Range("B12").Value = 16.06261469055171
$ws.Range("C12").Value = 8.047607746551181
$ws.Range("E12").Value = 17.63504386510577
$ws.Range("F12").Value = 49.50948221170858
$ws.Range("G12").Value = 3.696241597550211
$ws.Range("J12").Value = 9.566956366857898
$ws.Range("K12").Value = 15.16509240996774
$ws.Range("N12").Value = 21.66068471055528

# Row 13
$ws.Range("B13").Value = 16.0472259627333
$ws.Range("C13").Value = 8.036747715955464
$ws.Range("E13").Value = 17.62024298366814
$ws.Range("F13").Value = 49.48649033750085
$ws.Range("G13").Value = 3.696407713783852
$ws.Range("J13").Value = 9.567857051250177
$ws.Range("K13").Value = 15.15436070060096
$ws.Range("N13").Value = 21.6634335855769

# Row 14
$ws.Range("B14").Value = 15.99704278603059
$ws.Range("C14").Value = 8.001284522743724
$ws.Range("E14").Value = 17.57200710116289
$ws.Range("F14").Value = 49.41175879287574
$ws.Range("G14").Value = 3.696951836144792
$ws.Range("J14").Value = 9.570825476395514
$ws.Range("K14").Value = 15.11940647995726
$ws.Range("N14").Value = 21.67244807714514

# Row 15
$ws.Range("B15").Value = 15.96630756970517
$ws.Range("C15").Value = 7.979527665095502
$ws.Range("E15").Value = 17.54248777045646
$ws.Range("F15").Value = 49.36617694905058
$ws.Range("G15").Value = 3.697286958104371
$ws.Range("J15").Value = 9.572667615292557
$ws.Range("K15").Value = 15.09803077774468
$ws.Range("N15").Value = 21.67800798587491

# Row 16
$ws.Range("B16").Value = 15.79028185657905
$ws.Range("C16").Value = 7.85435624758926
$ws.Range("E16").Value = 17.37377994452067
$ws.Range("F16").Value = 49.1079836584661
$ws.Range("G16").Value = 3.69923553880845
$ws.Range("J16").Value = 9.583590248513323
$ws.Range("K16").Value = 14.97610081966448
$ws.Range("N16").Value = 21.71045600465568

# Row 17
$ws.Range("B17").Value = 15.6824478638517
$ws.Range("C17").Value = 7.777172490790919
$ws.Range("E17").Value = 17.27074289327329
$ws.Range("F17").Value = 48.95234088921239
$ws.Range("G17").Value = 3.700456101825858
$ws.Range("J17").Value = 9.590617446139671
$ws.Range("K17").Value = 14.90184090386275
$ws.Range("N17").Value = 21.7308849059786

# Row 18
$ws.Range("B18").Value = 15.62050373298355
$ws.Range("C18").Value = 7.732650937563118
$ws.Range("E18").Value = 17.21166886867264
$ws.Range("F18").Value = 48.86385251655641
$ws.Range("G18").Value = 3.701167415357232
$ws.Range("J18").Value = 9.594779240464604
$ws.Range("K18").Value = 14.85934113996371
$ws.Range("N18").Value = 21.74282723319411

# Row 19
$ws.Range("B19").Value = 15.59954632909118
$ws.Range("C19").Value = 7.717556528039291
$ws.Range("E19").Value = 17.19170209825399
$ws.Range("F19").Value = 48.83407142405504
$ws.Range("G19").Value = 3.701409849995738
$ws.Range("J19").Value = 9.596208946162541
$ws.Range("K19").Value = 14.84498929198124
$ws.Range("N19").Value = 21.7469036977641

# Row 20
$ws.Range("B20").Value = 15.69391931445266
$ws.Range("C20").Value = 7.785402435649649
$ws.Range("E20").Value = 17.28169216178927
$ws.Range("F20").Value = 48.96880289223199
$ws.Range("G20").Value = 3.700325211196452
$ws.Range("J20").Value = 9.589856974744757
$ws.Range("K20").Value = 14.90972431480478
$ws.Range("N20").Value = 21.72869032494211

# Row 21
$ws.Range("B21").Value = 16.01178221266477
$ws.Range("C21").Value = 8.01170823078615
$ws.Range("E21").Value = 17.58616973444091
$ws.Range("F21").Value = 49.43366908754315
$ws.Range("G21").Value = 3.69679163272694
$ws.Range("J21").Value = 9.569948598425642
$ws.Range("K21").Value = 15.12966621456193
$ws.Range("N21").Value = 21.66979232607921

# Row 22
$ws.Range("B22").Value = 16.21973965953723
$ws.Range("C22").Value = 8.158109241878883
$ws.Range("E22").Value = 17.78640828715949
$ws.Range("F22").Value = 49.74620259359973
$ws.Range("G22").Value = 3.694564273300765
$ws.Range("J22").Value = 9.55800655406431
$ws.Range("K22").Value = 15.27500497768703
$ws.Range("N22").Value = 21.63301186153644

# Row 23
$ws.Range("B23").Value = 16.10875238133066
$ws.Range("C23").Value = 8.080126845098823
$ws.Range("E23").Value = 17.6794448159574
$ws.Range("F23").Value = 49.57862402758437
$ws.Range("G23").Value = 3.695745578465374
$ws.Range("J23").Value = 9.564282355009766
$ws.Range("K23").Value = 15.19730355612236
$ws.Range("N23").Value = 21.65248548502386

# Row 24
$ws.Range("B24").Value = 15.68873290704608
$ws.Range("C24").Value = 7.781682133084599
$ws.Range("E24").Value = 17.27674148369393
$ws.Range("F24").Value = 48.96135732305015
$ws.Range("G24").Value = 3.70038435695401
$ws.Range("J24").Value = 9.590200404595556
$ws.Range("K24").Value = 14.90615961942163
$ws.Range("N24").Value = 21.72968187987348

# Row 25
$ws.Range("B25").Value = 15.23900382926386
$ws.Range("C25").Value = 7.455138540409362
$ws.Range("E25").Value = 16.84990081831638
$ws.Range("F25").Value = 48.33528906675037
$ws.Range("G25").Value = 3.705743106325436
$ws.Range("J25").Value = 9.622741884273717
$ws.Range("K25").Value = 14.60041455864276
$ws.Range("N25").Value = 21.82029699492069

Write-Host "Applied 192 cell updates"
